# aggiornamento fino a 13/03
# Adds 4 new daily-data rows (252-255) below the existing last row (251),
# following the same layout: column A is the date serial, columns B..AX
# are the per-comune counts / total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 251
$newRowCount = 4

# New data, one array per new row: Date, then values for columns B..AX (49 values)
$newData = @(
    @(44326,17,4,0,66,43,8,15,8,5,2,8,17,21,0,1,0,26,5,10,16,211,8,21,18,50,3,0,2,6,7,2,56,10,8,2,10,6,27,5,20,752,5,0,1,0,0,1,1,0),
    @(44327,18,3,0,72,40,7,16,8,6,3,8,17,22,0,1,0,30,6,11,15,209,8,23,19,60,1,0,2,5,6,1,55,9,9,1,11,7,29,5,20,770,5,0,0,0,0,1,1,0),
    @(44328,18,3,0,71,41,6,14,7,6,3,8,14,24,0,1,1,30,6,12,15,211,8,23,19,61,1,1,3,4,6,1,49,9,8,1,9,8,28,5,17,759,5,0,0,0,0,1,1,0),
    @(44329,13,1,0,76,26,10,14,4,3,2,9,8,24,0,1,1,23,4,8,16,223,6,20,16,55,0,2,3,3,5,1,41,5,4,1,8,9,25,4,13,695,5,0,0,0,0,1,2,0)
)

for ($i = 0; $i -lt $newRowCount; $i++) {
    $targetRow = $lastExistingRow + 1 + $i

    # Copy formatting (number format, borders, alignment, font, etc.) from the
    # last existing data row down onto the new row first.
    $srcRange = $ws.Range("A" + $lastExistingRow + ":AX" + $lastExistingRow)
    $dstRange = $ws.Range("A" + $targetRow + ":AX" + $targetRow)
    $srcRange.Copy($dstRange)

    $rowValues = $newData[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($targetRow, $col).Value2 = $rowValues[$col - 1]
    }
}
